# Re-pull data / push all data / mean calculation
# Updates the dSF column (F) values on the active sheet to match the
# freshly re-pulled dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = -5
    4  = -8
    5  = 7
    6  = -1
    10 = -2
    11 = 4
    12 = -5
    13 = 2
    14 = 2
    15 = 5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
